$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 597.75
$ws.Range("I6").Value = 227.1
$ws.Range("J6").Value = 2451
$ws.Range("K6").Value = 681.3
$ws.Range("L6").Value = 7353
$ws.Range("M6").Value = -569.3
$ws.Range("N6").Value = -7577
$ws.Range("H33").Value = 422
$ws.Range("I33").Value = 261.5
$ws.Range("K33").Value = 261.5
$ws.Range("M33").Value = -32.5
$ws.Range("H38").Value = 840.3333
$ws.Range("J38").Value = 1653
$ws.Range("L38").Value = 4959
$ws.Range("N38").Value = -5703
$ws.Range("H39").Value = 614
$ws.Range("H41").Value = 1376.4615
$ws.Range("I41").Value = 329.5
$ws.Range("K41").Value = 329.5
$ws.Range("M41").Value = 110.5
$ws.Range("H64").Value = 7301.467
$ws.Range("I64").Value = 6891.4443
$ws.Range("K64").Value = 6891.4443
$ws.Range("M64").Value = -6643.4443
$ws.Range("H67").Value = 7301.467
$ws.Range("I67").Value = 6891.4443
$ws.Range("K67").Value = 6891.4443
$ws.Range("M67").Value = -6033.4443
$ws.Range("H100").Value = 10625
$ws.Range("I100").Value = 10625
$ws.Range("K100").Value = 10625
$ws.Range("M100").Value = -10084
$ws.Range("H116").Value = 21895.479
$ws.Range("I116").Value = 22435.53
$ws.Range("J116").Value = 20365.334
$ws.Range("K116").Value = 22435.53
$ws.Range("L116").Value = 20365.334
$ws.Range("M116").Value = -18993.53
$ws.Range("N116").Value = -27249.334
$ws.Range("H138").Value = 50628.477
$ws.Range("I138").Value = 3034.875
$ws.Range("J138").Value = 79916.84
$ws.Range("K138").Value = 9104.625
$ws.Range("L138").Value = 239750.52
$ws.Range("M138").Value = -3964.625
$ws.Range("N138").Value = -250030.52

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1262.4
$ws.Range("I2").Value = 2500
$ws.Range("J2").Value = 953
$ws.Range("K2").Value = 2500
$ws.Range("L2").Value = 953
$ws.Range("M2").Value = -2387
$ws.Range("N2").Value = -1179
$ws.Range("H63").Value = 2682.9473
$ws.Range("I63").Value = 2520.75
$ws.Range("K63").Value = 2520.75
$ws.Range("M63").Value = -1834.75
$ws.Range("H66").Value = 2682.9473
$ws.Range("I66").Value = 2520.75
$ws.Range("K66").Value = 12603.75
$ws.Range("M66").Value = -9171.75
$ws.Range("H97").Value = 1257.08
$ws.Range("I97").Value = 1053.5714
$ws.Range("J97").Value = 2325.5
$ws.Range("K97").Value = 1053.5714
$ws.Range("L97").Value = 2325.5
$ws.Range("M97").Value = -557.5714
$ws.Range("N97").Value = -3317.5
$ws.Range("H116").Value = 1262.4
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 953
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 953
$ws.Range("M116").Value = -206
$ws.Range("N116").Value = -5541
$ws.Range("H132").Value = 1406.0571
$ws.Range("J132").Value = 2821.7144
$ws.Range("L132").Value = 8465.143199999999
$ws.Range("N132").Value = -13525.1432
$ws.Range("H135").Value = 64714.5
$ws.Range("J135").Value = 64714.5
$ws.Range("L135").Value = 64714.5
$ws.Range("N135").Value = -74854.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1262.4
$ws.Range("I3").Value = 2500
$ws.Range("J3").Value = 953
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 953
$ws.Range("M3").Value = -2386
$ws.Range("N3").Value = -1181
$ws.Range("H105").Value = 1817.9656
$ws.Range("I105").Value = 1215.8235
$ws.Range("J105").Value = 2671
$ws.Range("K105").Value = 1215.8235
$ws.Range("L105").Value = 2671
$ws.Range("M105").Value = 531.1765
$ws.Range("N105").Value = -6165

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3145
$ws.Range("I16").Value = 1784.5
$ws.Range("J16").Value = 4505.5
$ws.Range("K16").Value = 1784.5
$ws.Range("L16").Value = 4505.5
$ws.Range("M16").Value = -1497.5
$ws.Range("N16").Value = -5079.5
$ws.Range("H113").Value = 3145
$ws.Range("I113").Value = 1784.5
$ws.Range("J113").Value = 4505.5
$ws.Range("K113").Value = 1784.5
$ws.Range("L113").Value = 4505.5
$ws.Range("M113").Value = 385.5
$ws.Range("N113").Value = -8845.5
$ws.Range("H134").Value = 2718.0688
$ws.Range("I134").Value = 2564.0417
$ws.Range("J134").Value = 3457.4
$ws.Range("K134").Value = 7692.125100000001
$ws.Range("L134").Value = 10372.2
$ws.Range("M134").Value = -5157.125100000001
$ws.Range("N134").Value = -15442.2

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 623.8095
$ws.Range("I2").Value = 891.8182
$ws.Range("J2").Value = 329
$ws.Range("K2").Value = 891.8182
$ws.Range("L2").Value = 329
$ws.Range("M2").Value = -778.8182
$ws.Range("N2").Value = -555
$ws.Range("H97").Value = 1206.3636
$ws.Range("I97").Value = 1206.3636
$ws.Range("K97").Value = 1206.3636
$ws.Range("M97").Value = -710.3635999999999
$ws.Range("H113").Value = 1490.8125
$ws.Range("I113").Value = 1490.8125
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1490.8125
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 679.1875
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2991.6
$ws.Range("I122").Value = 2939.1765
$ws.Range("J122").Value = 3103
$ws.Range("K122").Value = 8817.529500000001
$ws.Range("L122").Value = 9309
$ws.Range("M122").Value = -6367.529500000001
$ws.Range("N122").Value = -14209
$ws.Range("H132").Value = 1969.25
$ws.Range("I132").Value = 1873.5278
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 5620.5834
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -3090.5834
$ws.Range("N132").Value = -12260

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7750
$ws.Range("I7").Value = 5333.3335
$ws.Range("J7").Value = 15000
$ws.Range("K7").Value = 5333.3335
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = -5221.3335
$ws.Range("N7").Value = -15224
$ws.Range("H16").Value = 1199.0834
$ws.Range("I16").Value = 1199.0834
$ws.Range("K16").Value = 1199.0834
$ws.Range("M16").Value = -1029.0834
$ws.Range("H46").Value = 2495.5518
$ws.Range("I46").Value = 1236.5454
$ws.Range("K46").Value = 1236.5454
$ws.Range("M46").Value = -1048.5454
$ws.Range("H61").Value = 571
$ws.Range("I61").Value = 571
$ws.Range("K61").Value = 571
$ws.Range("M61").Value = -369
$ws.Range("H113").Value = 571
$ws.Range("I113").Value = 571
$ws.Range("K113").Value = 571
$ws.Range("M113").Value = 1599
$ws.Range("H122").Value = 2251.4211
$ws.Range("I122").Value = 2210.4119
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 6631.2357
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -4181.2357
$ws.Range("N122").Value = -12700
$ws.Range("H126").Value = 7750
$ws.Range("I126").Value = 5333.3335
$ws.Range("J126").Value = 15000
$ws.Range("K126").Value = 16000.0005
$ws.Range("L126").Value = 45000
$ws.Range("M126").Value = -13530.0005
$ws.Range("N126").Value = -49940
$ws.Range("H136").Value = 2826.4878
$ws.Range("I136").Value = 2194.5806
$ws.Range("J136").Value = 4785.4
$ws.Range("K136").Value = 6583.7418
$ws.Range("L136").Value = 14356.2
$ws.Range("M136").Value = -4033.7418
$ws.Range("N136").Value = -19456.2

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 29659.666
$ws.Range("J104").Value = 29659.666
$ws.Range("L104").Value = 29659.666
$ws.Range("N104").Value = -36647.666
$ws.Range("H132").Value = 20932.629
$ws.Range("I132").Value = 25654.428
$ws.Range("J132").Value = 2045.4286
$ws.Range("K132").Value = 76963.284
$ws.Range("L132").Value = 6136.2858
$ws.Range("M132").Value = -74433.284
$ws.Range("N132").Value = -11196.2858
$ws.Range("H137").Value = 109742.8
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 109742.8
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 109742.8
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -119942.8
